$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Mona", "", "0555935549", "", "Polestar 2", "Saturday", "", "", "", "test drive"),
    @("عبدالله محمد عبدالرحيم", "", "65444356694", "76905", "تايوتا هايليندر مودل عام ٢٠٢٣", "2025-05-01", "", "", "12:00", ""),
    @("عبدالله محمد عبدالرهيب", "", "0500556694", "76905", "تيوتا هايلندر 2023", "1-5-2025", "", "", "", ""),
    @("Mona", "", "0555935549", "", "Polestar 2", "Saturday", "", "", "", "test drive"),
    @("Mohammed Al-Naim", "", "", "", "Dodge RAM 1500 Limited", "", "", "", "", "")
)

$startRow = 17
$endRow = $startRow + $data.Length - 1

# Force all the new cells to be formatted as Text so that numeric-looking
# strings (phone numbers, date-like strings) are preserved exactly as text
# and are not auto-converted to numbers or dates.
$ws.Range("A$startRow`:J$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 10; $col++) {
        $val = $rowData[$col - 1]
        $ws.Cells.Item($row, $col).Value = $val
    }
}
